$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column P into the new column Q for rows 3-34,
# then set the 2020 data values (or "-" placeholder) in Q4:Q34.
$ws.Range("P3").Copy($ws.Range("Q3"))

$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("Q4").Value = 2020

$ws.Range("P5").Copy($ws.Range("Q5"))
$ws.Range("Q5").Value = 51

$ws.Range("P6").Copy($ws.Range("Q6"))
$ws.Range("Q6").Value = 29

$ws.Range("P7").Copy($ws.Range("Q7"))
$ws.Range("Q7").Value = 22

$ws.Range("P8").Copy($ws.Range("Q8"))
$ws.Range("Q8").Value = 5

$ws.Range("P9").Copy($ws.Range("Q9"))
$ws.Range("Q9").Value = 3

$ws.Range("P10").Copy($ws.Range("Q10"))
$ws.Range("Q10").Value = 2

$ws.Range("P11").Copy($ws.Range("Q11"))
$ws.Range("Q11").Value = 15

$ws.Range("P12").Copy($ws.Range("Q12"))
$ws.Range("Q12").Value = 9

$ws.Range("P13").Copy($ws.Range("Q13"))
$ws.Range("Q13").Value = 5

$ws.Range("P14").Copy($ws.Range("Q14"))
$ws.Range("Q14").Value = "-"

$ws.Range("P15").Copy($ws.Range("Q15"))
$ws.Range("Q15").Value = "-"

$ws.Range("P16").Copy($ws.Range("Q16"))
$ws.Range("Q16").Value = "-"

$ws.Range("P17").Copy($ws.Range("Q17"))
$ws.Range("Q17").Value = "-"

$ws.Range("P18").Copy($ws.Range("Q18"))
$ws.Range("Q18").Value = "-"

$ws.Range("P19").Copy($ws.Range("Q19"))
$ws.Range("Q19").Value = "-"

$ws.Range("P20").Copy($ws.Range("Q20"))
$ws.Range("Q20").Value = 7

$ws.Range("P21").Copy($ws.Range("Q21"))
$ws.Range("Q21").Value = 7

$ws.Range("P22").Copy($ws.Range("Q22"))
$ws.Range("Q22").Value = "-"

$ws.Range("P23").Copy($ws.Range("Q23"))
$ws.Range("Q23").Value = "-"

$ws.Range("P24").Copy($ws.Range("Q24"))
$ws.Range("Q24").Value = "-"

$ws.Range("P25").Copy($ws.Range("Q25"))
$ws.Range("Q25").Value = "-"

$ws.Range("P26").Copy($ws.Range("Q26"))
$ws.Range("Q26").Value = 24

$ws.Range("P27").Copy($ws.Range("Q27"))
$ws.Range("Q27").Value = 10

$ws.Range("P28").Copy($ws.Range("Q28"))
$ws.Range("Q28").Value = 14

$ws.Range("P29").Copy($ws.Range("Q29"))
$ws.Range("Q29").Value = "-"

$ws.Range("P30").Copy($ws.Range("Q30"))
$ws.Range("Q30").Value = "-"

$ws.Range("P31").Copy($ws.Range("Q31"))
$ws.Range("Q31").Value = "-"

$ws.Range("P32").Copy($ws.Range("Q32"))
$ws.Range("Q32").Value = "-"

$ws.Range("P33").Copy($ws.Range("Q33"))
$ws.Range("Q33").Value = "-"

$ws.Range("P34").Copy($ws.Range("Q34"))
$ws.Range("Q34").Value = "-"

# Restore the selection shown when the workbook was last saved.
$ws.Range("H26").Select()
